# Apply the crypto-price/volume refresh captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.412.08'

$ws.Range("D3").Value = '2.983.41'
$ws.Range("E3").Value = '  +1.26%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.48'
$ws.Range("E5").Value = '  +2.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.27'
$ws.Range("E6").Value = '  +3.86%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  +2.25%  '

$ws.Range("D9").Value = '2.975.70'
$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("E10").Value = '  +3.74%  '

$ws.Range("E11").Value = '  +12.12%  '

$ws.Range("E12").Value = '  +1.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("E13").Value = '  +3.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.73'
$ws.Range("E14").Value = '  +2.96%  '

$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").Value = '3.476.87'
$ws.Range("E16").Value = '  +1.40%  '

$ws.Range("E17").Value = '  +2.40%  '

$ws.Range("D18").Value = '2.978.50'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").Value = '59.418.40'
$ws.Range("E19").Value = '  +2.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.07'
$ws.Range("E20").Value = '  +4.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("E21").Value = '  +1.98%  '

$ws.Range("E22").Value = '  +3.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.35'
$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.92'
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.28%  '

$ws.Range("B27").Value = 'FirstDigitalUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  +9.68%  '

$ws.Range("E29").Value = '  +2.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.74'
$ws.Range("E30").Value = '  +3.70%  '

$ws.Range("E31").Value = '  +9.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.71'
$ws.Range("E32").Value = '  +0.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.21'
$ws.Range("E33").Value = '  +4.91%  '

$ws.Range("D34").Value = '0.0₃0768'
$ws.Range("E34").Value = '  +9.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.89'
$ws.Range("E35").Value = '  +3.72%  '

$ws.Range("E36").Value = '  +4.35%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.68'
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  +4.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '400.51'
$ws.Range("E41").Value = '  +5.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0351'
$ws.Range("E42").Value = '  +1.20%  '

$ws.Range("D43").Value = '2.749.11'
$ws.Range("E43").Value = '  +1.90%  '

$ws.Range("E44").Value = '  -2.79%  '

$ws.Range("E45").Value = '  +6.31%  '

$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '34.87'
$ws.Range("E47").Value = '  +18.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.51'
$ws.Range("E48").Value = '  -1.27%  '

$ws.Range("E49").Value = '  +3.03%  '

$ws.Range("E50").Value = '  +2.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.32'
$ws.Range("E51").Value = '  +1.86%  '

